$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "name" column (B) for the existing rows first (A1:A5 data,
# then A6:A12 data which will become A7:A13 once the new row is inserted).
$ws.Range("B1").Value = "André"
$ws.Range("B2").Value = "Maurício"
$ws.Range("B3").Value = "Lucas"
$ws.Range("B4").Value = "Marcelo"
$ws.Range("B5").Value = "Jamile"
$ws.Range("B6").Value = "Bryann"
$ws.Range("B7").Value = "Alex"
$ws.Range("B8").Value = "Natália"
$ws.Range("B9").Value = "Isabela"
$ws.Range("B10").Value = "Luiz"
$ws.Range("B11").Value = "Beatriz"
$ws.Range("B12").Value = "Sandra"

# Insert a brand-new record as row 6, pushing the old rows 6-12 down to 7-13.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = 1254
$ws.Range("B6").Value = "Anderson"

# Append extra (non-numeric) characters onto the first two phone/doc cells,
# turning them from numbers into text.
$ws.Range("A2").Value = "1234abcs"
$ws.Range("A1").Value = "(11) 99447-9393abcd"

# Resize the columns to fit the new, wider contents (A auto-fits to 19
# characters wide; B gets an explicit width close to column A's old width + 4).
$ws.Columns("A:A").AutoFit()
$ws.Columns("A:A").ColumnWidth = 18.1
$ws.Columns("B:B").ColumnWidth = 17.7

# Leave the active cell / selection on B1.
[void]$ws.Range("B1").Select()
